$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.073.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.638.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  +5.65%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.102.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "61.078.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000145"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.639.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "356.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.432"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0871"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.16%  "
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("E36").Value = "  +9.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.886"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "298.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.645"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.102"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0564"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0239"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.984.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
